$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.020.94"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.215.61"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'291.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.19%  "
$ws.Range("D6").Value = "'87.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.45%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.471"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("D10").Value = "'30.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.93%  "
$ws.Range("D11").Value = "'0.0785"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").Value = "'6.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("D15").Value = "2.564.45"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "'14.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "2.220.27"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "'0.729"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("D19").Value = "39.961.93"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "'11.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +11.77%  "
$ws.Range("D21").Value = "0.0₃0886"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "'5.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "'65.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").Value = "'234.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").Value = "'2.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("D27").Value = "'1.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.70%  "
$ws.Range("D28").Value = "'22.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("D29").Value = "'2.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("D30").Value = "'9.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'152.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "'32.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").Value = "'4.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.61%  "
$ws.Range("D35").Value = "'0.0719"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.70%  "
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("D37").Value = "'2.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.67%  "
$ws.Range("D38").Value = "'16.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("D40").Value = "'0.0999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("D41").Value = "'1.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.68%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.096.64"
$ws.Range("E42").Value = "  +8.59%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'3.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.42%  "
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0270"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("D46").Value = "'18.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.16%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'10.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.86%  "
$ws.Range("D48").Value = "'2.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("D49").Value = "2.440.33"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "'69.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("D51").Value = "'1.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.14%  "
